$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet")

# Rows 8, 11, 14: clear the AL column value (PPn) leaving the cell empty
$ws.Range("AL8").ClearContents()
$ws.Range("AL11").ClearContents()
$ws.Range("AL14").ClearContents()

# Row 15: correct Jumlah/Total/DPP values, clear AL15
$ws.Range("AH15").Value = 20
$ws.Range("AI15").Value = 300000
$ws.Range("AK15").Value = 300000
$ws.Range("AL15").ClearContents()

# Row 16: correct Jumlah/Total/DPP values, clear AL16
$ws.Range("AH16").Value = 500
$ws.Range("AI16").Value = 4000000
$ws.Range("AK16").Value = 4000000
$ws.Range("AL16").ClearContents()
